# Update the "Reason for discontinuation" summary table with revised
# percentages / counts for the 2021 Q2 quarterly report (second draft).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, preserving the run formatting.
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $newText
}

# Row 2 = "Stolen"
Set-CellText $table 2 4 "37 (29.8)"

# Row 3 = "Completed Study"
Set-CellText $table 3 4 "28 (22.6)"

# Row 4 = "Non-Compliance"
Set-CellText $table 4 4 "26 (21.0)"

# Row 5 = "Lost"
Set-CellText $table 5 4 "12 (9.7)"

# Row 6 = "Defective"
Set-CellText $table 6 4 "9 (7.3)"

# Row 7 = "Sold or Gifted"
Set-CellText $table 7 4 "9 (7.3)"

# Row 8 = "Returned to Staff"
Set-CellText $table 8 4 "2 (1.6)"
